$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right after "2021-Q4" and before "总计".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1Sheet = $wb.Worksheets.Add($null, $afterSheet)
$q1Sheet.Name = "2022-Q1"

# Header row (row 1, columns B..H) - same layout as the other quarterly sheets.
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Data row (row 2) - fund holding detail for 2022-Q1.
$q1Sheet.Range("A2").Value = 0

$q1Sheet.Range("B2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "510200"

$q1Sheet.Range("C2").NumberFormat = "@"
$q1Sheet.Range("C2").Value = "汇安上证证券ETF"

$q1Sheet.Range("D2").NumberFormat = "@"
$q1Sheet.Range("D2").Value = "0.74"

$q1Sheet.Range("E2").NumberFormat = "@"
$q1Sheet.Range("E2").Value = "97.70"

$q1Sheet.Range("F2").NumberFormat = "@"
$q1Sheet.Range("F2").Value = "2.88"

$q1Sheet.Range("G2").NumberFormat = "@"
$q1Sheet.Range("G2").Value = "0.0213"

$q1Sheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q1 and shift the existing quarterly rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# New row 2: only the index + date are known for the just-started quarter.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"

# Re-number the running index in column A for the rows that shifted down
# (their date/count/value text moved automatically with the row insert).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
